$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update time on row 2 first (matches shared-string insertion order observed in target)
$ws.Range("H2").Value = "'12:00"

# Row 3: update claim number (E3, keep trailing space) and claim date (G3) as text (quote-prefixed)
$ws.Range("E3").Value = "'11111003116 "
$ws.Range("G3").Value = "'05/04/2021"

# Row 2: update claim number (E2) and claim date (G2) as text (quote-prefixed)
$ws.Range("E2").Value = "'11111003102"
$ws.Range("G2").Value = "'26/10/2021"

# Update selection to match D7
$ws.Range("D7").Select()
